$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column D (inserted into shared strings right after "Department")
$ws.Range("D1").Value = "Salary($/month)"

# Salary values for rows 2-5
$ws.Range("D2").Value = 500
$ws.Range("D3").Value = 300
$ws.Range("D4").Value = 400
$ws.Range("D5").Value = 100

# Total row label (appended to shared strings last, after "Clerk")
$ws.Range("C6").Value = "Total"

# Total formula
$ws.Range("D6").Formula = "=SUM(D2:D5)"
